$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 corresponds to "ano" = 2025
# Update total_customers (C6), new_customers (E6), new_rate (G6), returning_rate (H6)
$ws.Range("C6").Value = 430
$ws.Range("E6").Value = 121
$ws.Range("G6").Value = 28.13953488372093
$ws.Range("H6").Value = 71.86046511627907
